$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.594.17"
$ws.Range("E2").Value = "  +2.38%  "

$ws.Range("D3").Value = "1.792.75"
$ws.Range("E3").Value = "  +4.31%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5338"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +9.74%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3777"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.79%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.78"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07543"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.97%  "

$ws.Range("E11").Value = "  +6.54%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9993"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.09%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.190"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.72%  "

$ws.Range("D15").Value = "1.789.40"
$ws.Range("E15").Value = "  +4.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.095"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001069"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.96%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06508"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9990"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.12%  "

$ws.Range("E21").Value = "  +2.87%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.944"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.32%  "

$ws.Range("D23").Value = "27.629.81"
$ws.Range("E23").Value = "  +2.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.092"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.54%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.53"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.91%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.403"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +16.09%  "

$ws.Range("D29").Value = "1.994.29"
$ws.Range("E29").Value = "  +4.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "122.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.120"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1030"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +10.90%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.701"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.612"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.77%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02288"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.94%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.642"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +14.99%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06039"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.73%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.008"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.83%  "

$ws.Range("E39").Value = "  +4.32%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2088"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.80%  "

$ws.Range("E41").Value = "  +3.99%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.412"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.45%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9989"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.148"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.67%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5883"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.44%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.642"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.92%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "121.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.922"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.134"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06750"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.51%  "

